$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 201 in column C currently hold 7310; update them to 7293.
$ws.Range("C2:C201").Value = 7293
